$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Insert a new column before C (shifts C:N -> D:O), duplicating formatting
# from the old column C the way Excel's own "Insert Column" command does.
$ws.Columns("C:C").Insert()

# The insert also slides the "Date Range" header block (old C1:J1 / C2:F2)
# one column to the right, which is not what the human edit did there - in
# the source workbook that block stays anchored at C/D/E/F. Move it back.
$ws.Range("D1").Value = $ws.Range("C1").Value
$ws.Range("C1").Value = "Apiary Hive Inspection Regional Summary Report"
$ws.Range("D1").ClearContents()

$ws.Range("C2").Value = "Date Range:"
$ws.Range("D2").Value = "{d.DateRangeStart}"
$ws.Range("E2").Value = "TO"
$ws.Range("F2").Value = "{d.DateRangeEnd}"
$ws.Range("G2").ClearContents()

# New "Colonies Inspected Brood" column the insert created at C6:C13.
$ws.Range("C6").Value = "Colonies Inspected Brood"
$ws.Range("C7").Value = "{d.Region[i].BroodsInspected}"
$ws.Range("C13").Value = "{d.Tot_Broods_Inspected}"

# The Nosema running-total placeholder is renamed (and now sits at F13
# after the column insert shifted the totals row).
$ws.Range("F13").Value = "{d.Tot_Nosema}"

# New column C gets an explicit width (matches column B's width, 25,
# without "best fit").
$ws.Columns("C:C").ColumnWidth = 25

# Sheet view: selection moves to C14, and the frozen/top-left anchor on
# G1 is dropped.
$ws.Range("C14").Select()

$wb.Save()
